$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 265
$ws.Range("F4").Value = 592
$ws.Range("F5").Value = 2601
$ws.Range("F7").Value = 182
$ws.Range("F9").Value = 255
$ws.Range("F10").Value = 5554
$ws.Range("F11").Value = 108
$ws.Range("F12").Value = 1484
$ws.Range("F13").Value = 1417
$ws.Range("F14").Value = 610
$ws.Range("F16").Value = 52
$ws.Range("F19").Value = 4734
$ws.Range("F21").Value = 75
$ws.Range("F22").Value = 2380
$ws.Range("F23").Value = 1269
$ws.Range("F24").Value = 451
$ws.Range("F25").Value = 1170
$ws.Range("F26").Value = 227
$ws.Range("F27").Value = 94
$ws.Range("F28").Value = 80
$ws.Range("F29").Value = 169
$ws.Range("F30").Value = 371
$ws.Range("F31").Value = 1289
$ws.Range("F32").Value = 1999
$ws.Range("F33").Value = 246
$ws.Range("F34").Value = 532
$ws.Range("F35").Value = 12
$ws.Range("F36").Value = 207
$ws.Range("F37").Value = 1382
$ws.Range("F40").Value = 524
$ws.Range("F41").Value = 183
$ws.Range("F42").Value = 1633
$ws.Range("F43").Value = 2430
$ws.Range("F47").Value = 237
$ws.Range("F48").Value = 80
$ws.Range("F49").Value = 24

$ws = $wb.Worksheets.Item(2)
$ws.Range("F13").Value = 274
$ws.Range("F15").Value = 49
$ws.Range("F16").Value = 187
$ws.Range("F20").Value = 136
$ws.Range("F23").Value = 141
$ws.Range("F27").Value = 307
$ws.Range("F28").Value = 292

$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 500
$ws.Range("F5").Value = 19
$ws.Range("F6").Value = 1665
$ws.Range("F7").Value = 539
$ws.Range("F8").Value = 1326
$ws.Range("F9").Value = 1202
$ws.Range("F10").Value = 1760
$ws.Range("F11").Value = 2235
$ws.Range("F12").Value = 679
$ws.Range("F13").Value = 552

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 266
$ws.Range("F3").Value = 539
$ws.Range("F4").Value = 2601
$ws.Range("F5").Value = 182
$ws.Range("F6").Value = 1326
$ws.Range("F7").Value = 255
$ws.Range("F8").Value = 2235
$ws.Range("F9").Value = 5554
$ws.Range("F10").Value = 679
$ws.Range("F14").Value = 108
$ws.Range("F16").Value = 1484
$ws.Range("F17").Value = 1417
$ws.Range("F21").Value = 52
$ws.Range("F22").Value = 4735
$ws.Range("F23").Value = 2380
$ws.Range("F24").Value = 1170
$ws.Range("F25").Value = 94
$ws.Range("F26").Value = 80
$ws.Range("F27").Value = 274
$ws.Range("F29").Value = 169
$ws.Range("F30").Value = 49
$ws.Range("F31").Value = 187
$ws.Range("F32").Value = 371
$ws.Range("F33").Value = 1289
$ws.Range("F34").Value = 1999
$ws.Range("F35").Value = 532
$ws.Range("F36").Value = 136
$ws.Range("F37").Value = 207
$ws.Range("F38").Value = 1382
$ws.Range("F39").Value = 141
$ws.Range("F42").Value = 307
$ws.Range("F44").Value = 1633
$ws.Range("F45").Value = 2430
$ws.Range("F47").Value = 237
$ws.Range("F48").Value = 80
$ws.Range("F49").Value = 24
